$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 763.3333
$ws.Range("I28").Value = 443.55554
$ws.Range("J28").Value = 1083.1111
$ws.Range("K28").Value = 443.55554
$ws.Range("L28").Value = 1083.1111
$ws.Range("M28").Value = 41.44445999999999
$ws.Range("N28").Value = -2053.1111

$ws.Range("H111").Value = 8344098
$ws.Range("I111").Value = 26914.5
$ws.Range("J111").Value = 12502690
$ws.Range("K111").Value = 80743.5
$ws.Range("L111").Value = 37508070
$ws.Range("M111").Value = -77676.5
$ws.Range("N111").Value = -37514204

$ws.Range("H112").Value = 1147.826
$ws.Range("J112").Value = 1157.7778
$ws.Range("L112").Value = 3473.3334
$ws.Range("N112").Value = -5689.3334

$ws.Range("H129").Value = 937.94446
$ws.Range("I129").Value = 647.3333
$ws.Range("J129").Value = 1083.25
$ws.Range("K129").Value = 1941.9999
$ws.Range("L129").Value = 3249.75
$ws.Range("M129").Value = 3058.0001
$ws.Range("N129").Value = -13249.75

$ws.Range("H132").Value = 5438361.5
$ws.Range("I132").Value = 5559125
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 16677375
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -16674845
$ws.Range("N132").Value = -17057

$ws.Range("H135").Value = 917
$ws.Range("I135").Value = 637.5333000000001
$ws.Range("J135").Value = 1965
$ws.Range("K135").Value = 5737.7997
$ws.Range("L135").Value = 17685
$ws.Range("M135").Value = -3202.7997
$ws.Range("N135").Value = -22755

$ws.Range("H137").Value = 1968.5358
$ws.Range("I137").Value = 1628.0454
$ws.Range("J137").Value = 3217
$ws.Range("K137").Value = 4884.1362
$ws.Range("L137").Value = 9651
$ws.Range("M137").Value = -2334.1362
$ws.Range("N137").Value = -14751

$ws.Range("H138").Value = 3575.592
$ws.Range("I138").Value = 1149.1621
$ws.Range("J138").Value = 5877.59
$ws.Range("K138").Value = 3447.4863
$ws.Range("L138").Value = 17632.77
$ws.Range("M138").Value = 1692.5137
$ws.Range("N138").Value = -27912.77

$ws.Range("H141").Value = 6630
$ws.Range("I141").Value = 6833.3335
$ws.Range("J141").Value = 6542.857
$ws.Range("K141").Value = 20500.0005
$ws.Range("L141").Value = 19628.571
$ws.Range("M141").Value = -15320.0005
$ws.Range("N141").Value = -29988.571


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 60484.41
$ws.Range("I2").Value = 1443.0769
$ws.Range("J2").Value = 252368.75
$ws.Range("K2").Value = 1443.0769
$ws.Range("L2").Value = 252368.75
$ws.Range("M2").Value = -1330.0769
$ws.Range("N2").Value = -252594.75

$ws.Range("H6").Value = 5833.6665
$ws.Range("I6").Value = 8002
$ws.Range("J6").Value = 5400
$ws.Range("K6").Value = 8002
$ws.Range("L6").Value = 5400
$ws.Range("M6").Value = -7829
$ws.Range("N6").Value = -5746

$ws.Range("H116").Value = 60484.41
$ws.Range("I116").Value = 1443.0769
$ws.Range("J116").Value = 252368.75
$ws.Range("K116").Value = 1443.0769
$ws.Range("L116").Value = 252368.75
$ws.Range("M116").Value = 850.9231
$ws.Range("N116").Value = -256956.75

$ws.Range("H135").Value = 43654.445
$ws.Range("J135").Value = 43654.445
$ws.Range("L135").Value = 43654.445
$ws.Range("N135").Value = -53794.445


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 60484.41
$ws.Range("I3").Value = 1443.0769
$ws.Range("J3").Value = 252368.75
$ws.Range("K3").Value = 1443.0769
$ws.Range("L3").Value = 252368.75
$ws.Range("M3").Value = -1329.0769
$ws.Range("N3").Value = -252596.75

$ws.Range("H94").Value = 59264.234
$ws.Range("I94").Value = 77376.62
$ws.Range("J94").Value = 399
$ws.Range("K94").Value = 77376.62
$ws.Range("L94").Value = 399
$ws.Range("M94").Value = -76925.62
$ws.Range("N94").Value = -1301

$ws.Range("H105").Value = 183676.19
$ws.Range("I105").Value = 127059.625
$ws.Range("K105").Value = 127059.625
$ws.Range("M105").Value = -125312.625

$ws.Range("H134").Value = 2537.4102
$ws.Range("I134").Value = 2433.6765
$ws.Range("K134").Value = 7301.029500000001
$ws.Range("M134").Value = -4766.029500000001


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1182.71
$ws.Range("I31").Value = 710.76544
$ws.Range("J31").Value = 3194.6843
$ws.Range("K31").Value = 710.76544
$ws.Range("L31").Value = 3194.6843
$ws.Range("M31").Value = -415.76544
$ws.Range("N31").Value = -3784.6843

$ws.Range("H34").Value = 1182.71
$ws.Range("I34").Value = 710.76544
$ws.Range("J34").Value = 3194.6843
$ws.Range("K34").Value = 710.76544
$ws.Range("L34").Value = 3194.6843
$ws.Range("M34").Value = -508.76544
$ws.Range("N34").Value = -3598.6843

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 4100.3794
$ws.Range("I107").Value = 9883.454
$ws.Range("J107").Value = 566.2778
$ws.Range("K107").Value = 9883.454
$ws.Range("L107").Value = 566.2778
$ws.Range("M107").Value = -7963.454
$ws.Range("N107").Value = -4406.2778


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6150.816
$ws.Range("I5").Value = 929.5
$ws.Range("J5").Value = 20770.5
$ws.Range("K5").Value = 2788.5
$ws.Range("L5").Value = 62311.5
$ws.Range("M5").Value = -2676.5
$ws.Range("N5").Value = -62535.5

$ws.Range("H125").Value = 1079.8
$ws.Range("I125").Value = 666.3333
$ws.Range("J125").Value = 1700
$ws.Range("K125").Value = 1998.9999
$ws.Range("L125").Value = 5100
$ws.Range("M125").Value = 2921.0001
$ws.Range("N125").Value = -14940

$ws.Range("H135").Value = 6150.816
$ws.Range("I135").Value = 929.5
$ws.Range("J135").Value = 20770.5
$ws.Range("K135").Value = 8365.5
$ws.Range("L135").Value = 186934.5
$ws.Range("M135").Value = -5830.5
$ws.Range("N135").Value = -192004.5


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2238.2307
$ws.Range("I132").Value = 1810.8422
$ws.Range("J132").Value = 3398.2856
$ws.Range("K132").Value = 5432.5266
$ws.Range("L132").Value = 10194.8568
$ws.Range("M132").Value = -2902.5266
$ws.Range("N132").Value = -15254.8568


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2164.9092
$ws.Range("I7").Value = 1337.3334
$ws.Range("J7").Value = 3158
$ws.Range("K7").Value = 1337.3334
$ws.Range("L7").Value = 3158
$ws.Range("M7").Value = -1225.3334
$ws.Range("N7").Value = -3382

$ws.Range("H40").Value = 85467.336
$ws.Range("I40").Value = 251402
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 251402
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -251266
$ws.Range("N40").Value = -2772

$ws.Range("H93").Value = 1185.909
$ws.Range("I93").Value = 1115.4615
$ws.Range("J93").Value = 1287.6666
$ws.Range("K93").Value = 1115.4615
$ws.Range("L93").Value = 1287.6666
$ws.Range("M93").Value = 132.5385000000001
$ws.Range("N93").Value = -3783.6666

$ws.Range("H126").Value = 2164.9092
$ws.Range("I126").Value = 1337.3334
$ws.Range("J126").Value = 3158
$ws.Range("K126").Value = 4012.0002
$ws.Range("L126").Value = 9474
$ws.Range("M126").Value = -1542.0002
$ws.Range("N126").Value = -14414


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1538.125
$ws.Range("J122").Value = 1826.25
$ws.Range("L122").Value = 5478.75
$ws.Range("N122").Value = -10378.75

$ws.Range("H126").Value = 2610
$ws.Range("I126").Value = 3460
$ws.Range("J126").Value = 1760
$ws.Range("K126").Value = 10380
$ws.Range("L126").Value = 5280
$ws.Range("M126").Value = -7910
$ws.Range("N126").Value = -10220

